# Apply the VLSP2016 column-1-analyze.xlsx update:
# Each row's column C holds an example list for the POS tag in column A.
# This change reorders/edits some of those example lists (shared strings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C4"  = ',, ., ", ..., :, -, ), (, ?, !'
    "C5"  = 'của, trong, với, cho, ở, để, về, từ, đến, trên'
    "C6"  = 'oxy, marketing, shop, nilông, cas, Moran, Oxy, radio, Games, vali'
    "C7"  = 'ơi, Ôi, ư, ạ, Ồ, Ừ, Vâng, nha, nhỉ, Hỡi'
    "C11" = 'người, chị, anh, bà, ông, cái, con, cô, chiếc, Anh'
    "C12" = 'Thuỳ, VN, Mỹ, Nguyễn, Trâm, Khiêm, Thành, miền, Văn, HCM'
    "C13" = 'đồng, phút, lít, g, giờ, ha, m, USD, tấn, giây'
    "C14" = 'TP., NTLS, ĐDV, BS, ĐĐV, CCV, TP, Q., AIDS, UBND'
    "C15" = 'mình, tôi, đó, này, ấy, họ, Tôi, đây, chúng tôi, gì'
    "C17" = 'cả, chính, rồi, thôi, ngay, thật, đâu, mà, nào, sao'
    "C19" = 'vì sao, một mình, như vậy, làm sao, Tại sao, Vì sao, ngày càng, vừa qua, thế nào, có lẽ'
    "C20" = 'đa, phó, phi, Phó, viên'
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
